$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55 (shifts existing rows 55..179 down to 56..180)
$ws.Rows("55:55").Insert()

# Populate the newly inserted row 55 with its data
$ws.Range("A55").Value = 4
$ws.Range("B55").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C55").Value = "Los Lagos"
$ws.Range("D55").Value = 44519
$ws.Range("E55").Value = 10
$ws.Range("F55").Value = 100112043
$ws.Range("G55").Value = "Pepino ensalada"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 11000
$ws.Range("L55").Value = 11000
$ws.Range("M55").Value = 11000
$ws.Range("N55").Value = "`$/caja 60 unidades"
$ws.Range("O55").Value = "Región de Arica y Parinacota"
$ws.Range("P55").Value = 183
$ws.Range("Q55").Value = 60
$ws.Range("R55").Value = "Hortaliza"
